$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates: volume number and report week dates
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# Column H best-fit width changed (content width changed)
$ws.Columns.Item(8).ColumnWidth = 6.71

# Crime data cell updates
$ws.Range("M14").Value = -69.230769230769
$ws.Range("N14").Value = -80
$ws.Range("F15").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("L15").Value = -43.75
$ws.Range("N15").Value = -83.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = -28.205128205128
$ws.Range("L16").Value = -26.956521739130
$ws.Range("M16").Value = -62.995594713656
$ws.Range("N16").Value = -89.756097560975
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = -5.263157894736
$ws.Range("I17").Value = 213
$ws.Range("J17").Value = 218
$ws.Range("K17").Value = -2.293577981651
$ws.Range("L17").Value = -7.792207792207
$ws.Range("M17").Value = -6.167400881057
$ws.Range("N17").Value = -66.350710900473
$ws.Range("C18").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -70
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -50.450450450450
$ws.Range("L18").Value = -58.955223880597
$ws.Range("M18").Value = -63.576158940397
$ws.Range("N18").Value = -89.889705882352
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = -18.181818181818
$ws.Range("I19").Value = 192
$ws.Range("J19").Value = 233
$ws.Range("K19").Value = -17.596566523605
$ws.Range("L19").Value = -32.394366197183
$ws.Range("M19").Value = -9.859154929577
$ws.Range("N19").Value = -24.110671936758
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 65
$ws.Range("J20").Value = 74
$ws.Range("K20").Value = -12.162162162162
$ws.Range("L20").Value = -32.291666666666
$ws.Range("M20").Value = 3.174603174603
$ws.Range("N20").Value = -84.560570071258
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -58.823529411764
$ws.Range("F21").Value = 58
$ws.Range("H21").Value = -17.142857142857
$ws.Range("I21").Value = 622
$ws.Range("J21").Value = 773
$ws.Range("K21").Value = -19.534282018111
$ws.Range("L21").Value = -29.398410896708
$ws.Range("M21").Value = -31.723380900109
$ws.Range("N21").Value = -77.340619307832
$ws.Range("C22").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 60
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = -15.789473684210
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 62
$ws.Range("K23").Value = 4.838709677419
$ws.Range("L23").Value = 12.068965517241
$ws.Range("M23").Value = 14.035087719298
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = -18.666666666666
$ws.Range("I24").Value = 480
$ws.Range("J24").Value = 593
$ws.Range("K24").Value = -19.055649241146
$ws.Range("L24").Value = -16.230366492146
$ws.Range("M24").Value = -12.087912087912
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 500
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 84
$ws.Range("J25").Value = 154
$ws.Range("K25").Value = -45.454545454545
$ws.Range("L25").Value = -51.724137931034
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = -11.111111111111
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 23.333333333333
$ws.Range("I26").Value = 296
$ws.Range("J26").Value = 353
$ws.Range("K26").Value = -16.147308781869
$ws.Range("L26").Value = 3.859649122807
$ws.Range("M26").Value = -47.048300536672
$ws.Range("F27").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = -45.454545454545
$ws.Range("D28").Value = 2
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -83.333333333333
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 40
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("M29").Value = -60
$ws.Range("N29").Value = -88.148148148148
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("M30").Value = -60
$ws.Range("N30").Value = -88.333333333333
